$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 - "Master Data" / "4) Book language" row
$ws.Range("E9").Value = "Language"
$ws.Range("G9").Value = "10 mins"
$ws.Range("I9").Value = "Completed"

# Row 10 - "Transaction Type" row
$ws.Range("E10").Value = "BookTran"
$ws.Range("G10").Value = "10 mins"
$ws.Range("I10").Value = "Completed"

# Update the active selection to I10, matching the author's final cursor position
$ws.Range("I10").Select()
